$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert two new columns at D:E, shifting existing D:K to F:M
$ws.Columns("D:E").Insert()

# Step 2: set values + number formats for the two new columns, row by row
$ws.Range("D7").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("E7").Value2 = 43373

$ws.Range("D8").NumberFormat = "#,##0"
$ws.Range("D8").Value2 = 425000
$ws.Range("E8").NumberFormat = "#,##0"
$ws.Range("E8").Value2 = 412000

$ws.Range("D9").NumberFormat = "#,##0"
$ws.Range("D9").Value2 = -64000
$ws.Range("E9").NumberFormat = "#,##0"
$ws.Range("E9").Value2 = -70000

$ws.Range("D10").NumberFormat = "#,##0"
$ws.Range("D10").Value2 = 489000
$ws.Range("E10").NumberFormat = "#,##0"
$ws.Range("E10").Value2 = 482000

$ws.Range("D11").NumberFormat = "#,##0"
$ws.Range("E11").NumberFormat = "#,##0"

$ws.Range("D12").NumberFormat = "#,##0"
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").NumberFormat = "#,##0"
$ws.Range("E12").Value2 = "NA"

$ws.Range("D13").NumberFormat = "#,##0"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").NumberFormat = "#,##0"
$ws.Range("E13").Value2 = 0

$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("D14").Value2 = 11000
$ws.Range("E14").NumberFormat = "#,##0"
$ws.Range("E14").Value2 = 24000

$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value2 = 0
$ws.Range("E15").NumberFormat = "#,##0"
$ws.Range("E15").Value2 = 0

$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("E16").NumberFormat = "#,##0"

$ws.Range("D17").NumberFormat = "#,##0"
$ws.Range("D17").Value2 = 318000
$ws.Range("E17").NumberFormat = "#,##0"
$ws.Range("E17").Value2 = 265000

$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value2 = 107000
$ws.Range("E18").NumberFormat = "#,##0"
$ws.Range("E18").Value2 = 147000

$ws.Range("D19").NumberFormat = "#,##0"
$ws.Range("E19").NumberFormat = "#,##0"

$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("D20").Value2 = 0
$ws.Range("E20").NumberFormat = "#,##0"
$ws.Range("E20").Value2 = 0

$ws.Range("D21").NumberFormat = "#,##0"
$ws.Range("D21").Value2 = "NA"
$ws.Range("E21").NumberFormat = "#,##0"
$ws.Range("E21").Value2 = "NA"

$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value2 = 0
$ws.Range("E22").NumberFormat = "#,##0"
$ws.Range("E22").Value2 = 0

$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("D23").Value2 = 107000
$ws.Range("E23").NumberFormat = "#,##0"
$ws.Range("E23").Value2 = 147000

$ws.Range("D24").NumberFormat = "#,##0"
$ws.Range("D24").Value2 = 35000
$ws.Range("E24").NumberFormat = "#,##0"
$ws.Range("E24").Value2 = 54000
$ws.Range("F24").Value2 = 46000
$ws.Range("G24").Value2 = 59000

$ws.Range("D25").NumberFormat = "#,##0"
$ws.Range("D25").Value2 = 0
$ws.Range("E25").NumberFormat = "#,##0"
$ws.Range("E25").Value2 = 0

$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value2 = 72000
$ws.Range("E26").NumberFormat = "#,##0"
$ws.Range("E26").Value2 = 93000
$ws.Range("F26").Value2 = 67000
$ws.Range("G26").Value2 = 103000

$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value2 = 72000
$ws.Range("E27").NumberFormat = "#,##0"
$ws.Range("E27").Value2 = 93000
$ws.Range("F27").Value2 = 67000
$ws.Range("G27").Value2 = 103000

$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value2 = 0
$ws.Range("E28").NumberFormat = "#,##0"
$ws.Range("E28").Value2 = 0

$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("D29").Value2 = 0
$ws.Range("E29").NumberFormat = "#,##0"
$ws.Range("E29").Value2 = 21000
$ws.Range("F29").Value2 = 16000
$ws.Range("G29").Value2 = 23000

$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Value2 = 0
$ws.Range("E30").NumberFormat = "#,##0"
$ws.Range("E30").Value2 = 0

$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("D31").Value2 = 0
$ws.Range("E31").NumberFormat = "#,##0"
$ws.Range("E31").Value2 = 0

$ws.Range("D32").NumberFormat = "#,##0"
$ws.Range("D32").Value2 = 0
$ws.Range("E32").NumberFormat = "#,##0"
$ws.Range("E32").Value2 = 0

$ws.Range("D33").NumberFormat = "#,##0"
$ws.Range("D33").Value2 = 72000
$ws.Range("E33").NumberFormat = "#,##0"
$ws.Range("E33").Value2 = 114000

$ws.Range("D34").NumberFormat = "#,##0"
$ws.Range("D34").Value2 = 0
$ws.Range("E34").NumberFormat = "#,##0"
$ws.Range("E34").Value2 = 0

$ws.Range("D35").NumberFormat = "#,##0"
$ws.Range("D35").Value2 = 72000
$ws.Range("E35").NumberFormat = "#,##0"
$ws.Range("E35").Value2 = 114000

$ws.Range("D38").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("E38").Value2 = 43373

$ws.Range("D39").NumberFormat = "#,##0"
$ws.Range("E39").NumberFormat = "#,##0"

$ws.Range("D40").NumberFormat = "#,##0"
$ws.Range("E40").NumberFormat = "#,##0"

$ws.Range("D41").NumberFormat = "#,##0"
$ws.Range("D41").Value2 = 1286000
$ws.Range("E41").NumberFormat = "#,##0"
$ws.Range("E41").Value2 = 2143000

$ws.Range("D42").NumberFormat = "#,##0"
$ws.Range("D42").Value2 = 0
$ws.Range("E42").NumberFormat = "#,##0"
$ws.Range("E42").Value2 = 0

$ws.Range("D43").NumberFormat = "#,##0"
$ws.Range("D43").Value2 = 2094000
$ws.Range("E43").NumberFormat = "#,##0"
$ws.Range("E43").Value2 = 2132000

$ws.Range("D44").NumberFormat = "#,##0"
$ws.Range("D44").Value2 = 0
$ws.Range("E44").NumberFormat = "#,##0"
$ws.Range("E44").Value2 = 0

$ws.Range("D45").NumberFormat = "#,##0"
$ws.Range("D45").Value2 = 0
$ws.Range("E45").NumberFormat = "#,##0"
$ws.Range("E45").Value2 = 0

$ws.Range("D46").NumberFormat = "#,##0"
$ws.Range("D46").Value2 = 0
$ws.Range("E46").NumberFormat = "#,##0"
$ws.Range("E46").Value2 = 0

$ws.Range("D47").NumberFormat = "#,##0"
$ws.Range("D47").Value2 = 95064000
$ws.Range("E47").NumberFormat = "#,##0"
$ws.Range("E47").Value2 = 97556000

$ws.Range("D48").NumberFormat = "#,##0"
$ws.Range("D48").Value2 = 136000
$ws.Range("E48").NumberFormat = "#,##0"
$ws.Range("E48").Value2 = 107000

$ws.Range("D49").NumberFormat = "#,##0"
$ws.Range("D49").Value2 = 786000
$ws.Range("E49").NumberFormat = "#,##0"
$ws.Range("E49").Value2 = 792000

$ws.Range("D50").NumberFormat = "#,##0"
$ws.Range("D50").Value2 = 0
$ws.Range("E50").NumberFormat = "#,##0"
$ws.Range("E50").Value2 = 0

$ws.Range("D51").NumberFormat = "#,##0"
$ws.Range("D51").Value2 = 0
$ws.Range("E51").NumberFormat = "#,##0"
$ws.Range("E51").Value2 = 0

$ws.Range("D52").NumberFormat = "#,##0"
$ws.Range("D52").Value2 = 4446000
$ws.Range("E52").NumberFormat = "#,##0"
$ws.Range("E52").Value2 = 3243000

$ws.Range("D53").NumberFormat = "#,##0"
$ws.Range("D53").Value2 = 0
$ws.Range("E53").NumberFormat = "#,##0"
$ws.Range("E53").Value2 = 0

$ws.Range("D54").NumberFormat = "#,##0"
$ws.Range("D54").Value2 = 104176000
$ws.Range("E54").NumberFormat = "#,##0"
$ws.Range("E54").Value2 = 106493000

$ws.Range("D55").NumberFormat = "#,##0"
$ws.Range("E55").NumberFormat = "#,##0"

$ws.Range("D56").NumberFormat = "#,##0"
$ws.Range("E56").NumberFormat = "#,##0"

$ws.Range("D57").NumberFormat = "#,##0"
$ws.Range("D57").Value2 = 0
$ws.Range("E57").NumberFormat = "#,##0"
$ws.Range("E57").Value2 = 0

$ws.Range("D58").NumberFormat = "#,##0"
$ws.Range("D58").Value2 = 5422000
$ws.Range("E58").NumberFormat = "#,##0"
$ws.Range("E58").Value2 = 5007000

$ws.Range("D59").NumberFormat = "#,##0"
$ws.Range("D59").Value2 = 0
$ws.Range("E59").NumberFormat = "#,##0"
$ws.Range("E59").Value2 = 0

$ws.Range("D60").NumberFormat = "#,##0"
$ws.Range("D60").Value2 = 0
$ws.Range("E60").NumberFormat = "#,##0"
$ws.Range("E60").Value2 = 0

$ws.Range("D61").NumberFormat = "#,##0"
$ws.Range("D61").Value2 = 93519000
$ws.Range("E61").NumberFormat = "#,##0"
$ws.Range("E61").Value2 = 96089000

$ws.Range("D62").NumberFormat = "#,##0"
$ws.Range("D62").Value2 = 0
$ws.Range("E62").NumberFormat = "#,##0"
$ws.Range("E62").Value2 = 0

$ws.Range("D63").NumberFormat = "#,##0"
$ws.Range("D63").Value2 = 0
$ws.Range("E63").NumberFormat = "#,##0"
$ws.Range("E63").Value2 = 0

$ws.Range("D64").NumberFormat = "#,##0"
$ws.Range("D64").Value2 = 0
$ws.Range("E64").NumberFormat = "#,##0"
$ws.Range("E64").Value2 = 0

$ws.Range("D65").NumberFormat = "#,##0"
$ws.Range("D65").Value2 = 0
$ws.Range("E65").NumberFormat = "#,##0"
$ws.Range("E65").Value2 = 0

$ws.Range("D66").NumberFormat = "#,##0"
$ws.Range("D66").Value2 = 100657000
$ws.Range("E66").NumberFormat = "#,##0"
$ws.Range("E66").Value2 = 102769000

$ws.Range("D67").NumberFormat = "#,##0"
$ws.Range("E67").NumberFormat = "#,##0"

$ws.Range("D68").NumberFormat = "#,##0"
$ws.Range("D68").Value2 = 0
$ws.Range("E68").NumberFormat = "#,##0"
$ws.Range("E68").Value2 = 0

$ws.Range("D69").NumberFormat = "#,##0"
$ws.Range("D69").Value2 = 0
$ws.Range("E69").NumberFormat = "#,##0"
$ws.Range("E69").Value2 = 0

$ws.Range("D70").NumberFormat = "#,##0"
$ws.Range("D70").Value2 = 0
$ws.Range("E70").NumberFormat = "#,##0"
$ws.Range("E70").Value2 = 0

$ws.Range("D71").NumberFormat = "#,##0"
$ws.Range("D71").Value2 = 0
$ws.Range("E71").NumberFormat = "#,##0"
$ws.Range("E71").Value2 = 0

$ws.Range("D72").NumberFormat = "#,##0"
$ws.Range("D72").Value2 = 3218000
$ws.Range("E72").NumberFormat = "#,##0"
$ws.Range("E72").Value2 = 3186000

$ws.Range("D73").NumberFormat = "#,##0"
$ws.Range("D73").Value2 = 0
$ws.Range("E73").NumberFormat = "#,##0"
$ws.Range("E73").Value2 = 0

$ws.Range("D74").NumberFormat = "#,##0"
$ws.Range("D74").Value2 = 0
$ws.Range("E74").NumberFormat = "#,##0"
$ws.Range("E74").Value2 = 0

$ws.Range("D75").NumberFormat = "#,##0"
$ws.Range("D75").Value2 = 0
$ws.Range("E75").NumberFormat = "#,##0"
$ws.Range("E75").Value2 = 0

$ws.Range("D76").NumberFormat = "#,##0"
$ws.Range("D76").Value2 = 3519000
$ws.Range("E76").NumberFormat = "#,##0"
$ws.Range("E76").Value2 = 3724000

$ws.Range("D77").NumberFormat = "#,##0"
$ws.Range("D77").Value2 = 0
$ws.Range("E77").NumberFormat = "#,##0"
$ws.Range("E77").Value2 = 0

$ws.Range("D80").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("E80").Value2 = 43373

$ws.Range("D81").NumberFormat = "#,##0"
$ws.Range("D81").Value2 = 72000
$ws.Range("E81").NumberFormat = "#,##0"
$ws.Range("E81").Value2 = 114000

$ws.Range("D82").NumberFormat = "#,##0"
$ws.Range("E82").NumberFormat = "#,##0"

$ws.Range("D83").NumberFormat = "#,##0"
$ws.Range("D83").Value2 = 0
$ws.Range("E83").NumberFormat = "#,##0"
$ws.Range("E83").Value2 = 0

$ws.Range("D84").NumberFormat = "#,##0"
$ws.Range("D84").Value2 = 0
$ws.Range("E84").NumberFormat = "#,##0"
$ws.Range("E84").Value2 = 0

$ws.Range("D85").NumberFormat = "#,##0"
$ws.Range("D85").Value2 = 0
$ws.Range("E85").NumberFormat = "#,##0"
$ws.Range("E85").Value2 = 0

$ws.Range("D86").NumberFormat = "#,##0"
$ws.Range("D86").Value2 = 0
$ws.Range("E86").NumberFormat = "#,##0"
$ws.Range("E86").Value2 = 0

$ws.Range("D87").NumberFormat = "#,##0"
$ws.Range("D87").Value2 = 0
$ws.Range("E87").NumberFormat = "#,##0"
$ws.Range("E87").Value2 = 0

$ws.Range("D88").NumberFormat = "#,##0"
$ws.Range("D88").Value2 = 0
$ws.Range("E88").NumberFormat = "#,##0"
$ws.Range("E88").Value2 = 0

$ws.Range("D89").NumberFormat = "#,##0"
$ws.Range("D89").Value2 = 485000
$ws.Range("E89").NumberFormat = "#,##0"
$ws.Range("E89").Value2 = 173000
$ws.Range("H89").Value2 = 244000

$ws.Range("D90").NumberFormat = "#,##0"
$ws.Range("E90").NumberFormat = "#,##0"

$ws.Range("D91").NumberFormat = "#,##0"
$ws.Range("D91").Value2 = 0
$ws.Range("E91").NumberFormat = "#,##0"
$ws.Range("E91").Value2 = 0

$ws.Range("D92").NumberFormat = "#,##0"
$ws.Range("D92").Value2 = 0
$ws.Range("E92").NumberFormat = "#,##0"
$ws.Range("E92").Value2 = 0

$ws.Range("D93").NumberFormat = "#,##0"
$ws.Range("D93").Value2 = 0
$ws.Range("E93").NumberFormat = "#,##0"
$ws.Range("E93").Value2 = 0

$ws.Range("D94").NumberFormat = "#,##0"
$ws.Range("D94").Value2 = 2161000
$ws.Range("E94").NumberFormat = "#,##0"
$ws.Range("E94").Value2 = 2443000

$ws.Range("D95").NumberFormat = "#,##0"
$ws.Range("E95").NumberFormat = "#,##0"

$ws.Range("D96").NumberFormat = "#,##0"
$ws.Range("D96").Value2 = -40000
$ws.Range("E96").NumberFormat = "#,##0"
$ws.Range("E96").Value2 = -41000

$ws.Range("D97").NumberFormat = "#,##0"
$ws.Range("D97").Value2 = 0
$ws.Range("E97").NumberFormat = "#,##0"
$ws.Range("E97").Value2 = 0

$ws.Range("D98").NumberFormat = "#,##0"
$ws.Range("D98").Value2 = 0
$ws.Range("E98").NumberFormat = "#,##0"
$ws.Range("E98").Value2 = 0

$ws.Range("D99").NumberFormat = "#,##0"
$ws.Range("D99").Value2 = 0
$ws.Range("E99").NumberFormat = "#,##0"
$ws.Range("E99").Value2 = 0

$ws.Range("D100").NumberFormat = "#,##0"
$ws.Range("D100").Value2 = -2632000
$ws.Range("E100").NumberFormat = "#,##0"
$ws.Range("E100").Value2 = -2376000
$ws.Range("H100").Value2 = -2411000
$ws.Range("I100").Value2 = -3330000

$ws.Range("D101").NumberFormat = "#,##0"
$ws.Range("D101").Value2 = 0
$ws.Range("E101").NumberFormat = "#,##0"
$ws.Range("E101").Value2 = 0

$ws.Range("D102").NumberFormat = "#,##0"
$ws.Range("D102").Value2 = 14000
$ws.Range("E102").NumberFormat = "#,##0"
$ws.Range("E102").Value2 = 240000
